$wb = $excel.ActiveWorkbook

# --- Sheet 1 (שבצ"כ) - rotation/shift corrections ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(3).ColumnWidth = 27.17

$ws1.Range("E88").Value = "עומרי דותן" + [char]10 + "דעאל כהן"
$ws1.Range("D89").Value = "אבנר יוזפוביץ" + [char]10 + "יניב משה"
$ws1.Range("C92").Value = "עדן אסרף" + [char]10 + "איתי סיני"
$ws1.Range("D92").Value = "ליאור אבו חמדה" + [char]10 + "מרדוש דהן"
$ws1.Range("E92").Value = "שראל בלוך" + [char]10 + "נתנאל שרעבי"
$ws1.Range("D98").Value = "אסף זבולון" + [char]10 + "אדיר מור"
$ws1.Range("C101").Value = "איתמר בנימין" + [char]10 + "מיכאל ניסנוב"
$ws1.Range("D101").Value = "דורון לביא" + [char]10 + "עמיחי נעים"
$ws1.Range("C110").Value = "לואיס אברבוך" + [char]10 + "זיטר יצחק"

# --- Sheet 2 (Available) - cumulative availability list corrections ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns.Item(3).ColumnWidth = 268.17

$ws2.Range("C44").Value = "שגיא אריה, דעאל כהן, עדן אסרף, איתי סיני, דובר אלבז" + [char]10 + "נריה כלפה, זיטר יצחק, נדב קריספין, יוסף רווה, ראובן מאור" + [char]10 + "ארד רז"
$ws2.Range("C47").Value = "דעאל כהן, אבנר יוזפוביץ, עדן אסרף, איתי סיני, דובר אלבז" + [char]10 + "נריה כלפה, חן טלה, ראובן מאור, יואל אודיז, ארד רז"
$ws2.Range("C48").Value = "דעאל כהן, אבנר יוזפוביץ, יניב משה, עדן אסרף, איתי סיני" + [char]10 + "דובר אלבז, נריה כלפה, חן טלה, אור נצקנסקי, ראובן מאור" + [char]10 + "יואל אודיז, ארד רז"
$ws2.Range("C76").Value = "אגומס מלדה, אבנר יוזפוביץ, ליאור אבו חמדה, יניב משה, עדן אסרף" + [char]10 + "איתי סיני, מרדוש דהן, שראל בלוך, נתנאל שרעבי, אייל דבוש" + [char]10 + "גיא פיאצה"
$ws2.Range("C77").Value = "ליאור אבו חמדה, עדן אסרף, איתי סיני, מרדוש דהן, שראל בלוך" + [char]10 + "נתנאל שרעבי, אייל דבוש, גיא פיאצה, דובר אלבז, נריה כלפה" + [char]10 + "חן טלה"
$ws2.Range("C80").Value = "אייל דבוש, גיא פיאצה, דובר אלבז, נריה כלפה, חן טלה" + [char]10 + "אסף זבולון, אדיר מור, מיכאל ניסנוב, לוטם עטיה"
$ws2.Range("C83").Value = "חן טלה, אסף זבולון, אדיר מור, מיכאל ניסנוב, לוטם עטיה" + [char]10 + "דורון לביא, עמיחי נעים, איתמר בנימין, זיטר יצחק"
$ws2.Range("C86").Value = "מיכאל ניסנוב, דורון לביא, עמיחי נעים, איתמר בנימין, זיטר יצחק" + [char]10 + "נדב קריספין, יוסף רווה, לואיס אברבוך, פביאן חויוס"
$ws2.Range("C89").Value = "זיטר יצחק, נדב קריספין, יוסף רווה, לואיס אברבוך, פביאן חויוס" + [char]10 + "אלכסיי ברומברג, ראובן מאור, יואל אודיז, ארד רז"
$ws2.Range("C92").Value = "דעאל כהן, אגומס מלדה, עומרי דותן, אבנר יוזפוביץ, יניב משה" + [char]10 + "זיטר יצחק, נדב קריספין, יוסף רווה, לואיס אברבוך, פביאן חויוס" + [char]10 + "אלכסיי ברומברג, ראובן מאור, יואל אודיז, ארד רז"
$ws2.Range("C95").Value = "דעאל כהן, אגומס מלדה, עומרי דותן, אבנר יוזפוביץ, ליאור אבו חמדה" + [char]10 + "יניב משה, עדן אסרף, איתי סיני, מרדוש דהן, זיטר יצחק" + [char]10 + "נדב קריספין, יוסף רווה, לואיס אברבוך, פביאן חויוס, אלכסיי ברומברג" + [char]10 + "ראובן מאור, יואל אודיז, ארד רז"
$ws2.Range("C96").Value = "שגיא אריה, דעאל כהן, אגומס מלדה, עומרי דותן, אבנר יוזפוביץ" + [char]10 + "ליאור אבו חמדה, יניב משה, עדן אסרף, איתי סיני, מרדוש דהן" + [char]10 + "שראל בלוך, נתנאל שרעבי, אור נצקנסקי, זיטר יצחק, נדב קריספין" + [char]10 + "יוסף רווה, לואיס אברבוך, פביאן חויוס, אלכסיי ברומברג, ראובן מאור" + [char]10 + "יואל אודיז, ארד רז"
$ws2.Range("C98").Value = "דעאל כהן, אגומס מלדה, עומרי דותן, אבנר יוזפוביץ, ליאור אבו חמדה" + [char]10 + "יניב משה, עדן אסרף, איתי סיני, מרדוש דהן, שראל בלוך" + [char]10 + "נתנאל שרעבי, אייל דבוש, גיא פיאצה, דובר אלבז, נריה כלפה" + [char]10 + "אור נצקנסקי, נדב קריספין, יוסף רווה, פביאן חויוס, אלכסיי ברומברג" + [char]10 + "ראובן מאור, יואל אודיז, ארד רז"
